# The edit cyclically rotates the "species record" data (columns A,B,D,E,F,G,H,
# the empty L marker cell, Q, R and the AC comment) among rows 2, 5, 7 and 9:
#   old row 2 data -> row 5
#   old row 5 data -> row 7
#   old row 7 data -> row 9
#   old row 9 data -> row 2
#
# Capture every source value up-front (using Value2, which this runtime's
# COM shim evaluates reliably) before any writes happen, so the rotation
# doesn't clobber values we still need to read.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2, 5, 7, 9)

$data = @{}
foreach ($r in $rows) {
    $row = @{
        A  = $ws.Range("A$r").Value2
        B  = $ws.Range("B$r").Value2
        D  = $ws.Range("D$r").Value2
        E  = $ws.Range("E$r").Value2
        F  = $ws.Range("F$r").Value2
        G  = $ws.Range("G$r").Value2
        H  = $ws.Range("H$r").Value2
        L  = $ws.Range("L$r").Value2
        Q  = $ws.Range("Q$r").Value2
        R  = $ws.Range("R$r").Value2
        AC = $ws.Range("AC$r").Value2
    }
    $data[$r] = $row
}

# Destination for each row's current content, per the rotation above.
$dest = @{ 2 = 5; 5 = 7; 7 = 9; 9 = 2 }

foreach ($src in $rows) {
    $target = $dest[$src]
    $vals = $data[$src]

    $ws.Range("A$target").Value = $vals.A
    $ws.Range("B$target").Value = $vals.B
    $ws.Range("D$target").Value = $vals.D
    $ws.Range("E$target").Value = $vals.E
    $ws.Range("F$target").Value = $vals.F
    $ws.Range("G$target").Value = $vals.G
    $ws.Range("H$target").Value = $vals.H

    if ($vals.L -eq $null) {
        $ws.Range("L$target").ClearContents()
    } else {
        $ws.Range("L$target").Value = $vals.L
    }

    $ws.Range("Q$target").Value = $vals.Q
    $ws.Range("R$target").Value = $vals.R

    if ($vals.AC -eq $null) {
        $ws.Range("AC$target").ClearContents()
    } else {
        $ws.Range("AC$target").Value = $vals.AC
    }
}
